# Apply the edit described in the commit:
# "SUMIF now copes with criteria ranges that are different in size from
#  sum ranges" — add a small demonstration block (rows 28-34) exercising
# SUMIF / SUMIFS with mismatched criteria/sum range shapes, and turn the
# row 19 formulas (C19:H19) into a proper shared-formula group.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 19: re-enter as one pass so Excel stores it as a shared formula
#     group (matches si="6" in the target workbook) instead of six
#     independent <f> elements.
$ws.Range("C19:H19").Formula = "=`$C`$6:`$H`$6"

# --- New demonstration data, columns E/F (small, 2 categories) and
#     H/I/J (mismatched-size criteria/sum ranges), rows 28-32.
$ws.Range("E28").Value = 1
$ws.Range("E29").Value = 2
$ws.Range("E30").Value = 3
$ws.Range("E31").Value = 4
$ws.Range("E32").Value = 5000

$ws.Range("F28").Value = "a"
$ws.Range("F29").Value = "b"
$ws.Range("F30").Value = "a"
$ws.Range("F31").Value = "b"

$ws.Range("H28").Value = 1
$ws.Range("H29").Value = 2
$ws.Range("H30").Value = 3000
$ws.Range("H31").Value = 4
$ws.Range("H32").Value = 5

$ws.Range("I28").Value = "a"
$ws.Range("I29").Value = "b"
$ws.Range("I30").Value = "a"
$ws.Range("I31").Value = "b"

$ws.Range("J28").Value = "x"
$ws.Range("J29").Value = "x"
$ws.Range("J30").Value = "x"
$ws.Range("J31").Value = "x"

# --- Formulas demonstrating SUMIF / SUMIFS with mismatched range sizes.
$ws.Range("E33").Formula = "=SUMIF(F28:F29,""a"",E28:E32)"
$ws.Range("H33").Formula = "=SUMIF(I28:I31,""a"",H28)"
$ws.Range("H34").Formula = "=SUMIFS(H28:H31,I28:I30,""a"",J28:J31,""x"")"

# --- Leave the selection where the author left it.
$ws.Range("H33").Select() | Out-Null
